{"js": "// Fix typo in BiDS'17 paper: \"Prolusion\" -> \"Propulsion\" in the first-page\n// affiliation footnote (\"* NASA Jet Prolusion Laboratory, California\n// Institute of Technology\"). Word also re-anchors its internal \"last edit\"\n// bookmark (_GoBack) to the location of the correction, so we mirror that\n// too: remove it from its old position and drop it right after the fix.\n\nconst doc = context.document;\n\n// 1) Locate and fix the misspelled word. It appears exactly once in the\n//    document (a correctly-spelled \"Jet Propulsion Laboratory (JPL)\" exists\n//    elsewhere and must be left untouched), so an exact, case-sensitive\n//    search is safe and unambiguous.\nconst misspelling = doc.body.search(\"Prolusion\", { matchCase: true, matchWholeWord: true });\nmisspelling.load(\"text\");\nawait context.sync();\n\nif (misspelling.items.length > 0) {\n  misspelling.items[0].insertText(\"Propulsion\", Word.InsertLocation.replace);\n  await context.sync();\n}\n\n// 2) Move the \"_GoBack\" bookmark (Word drops this at the point of the most\n//    recent edit) from wherever it currently sits to right after the word\n//    we just corrected.\ndoc.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst anchor = doc.body.search(\"NASA Jet Propul\", { matchCase: true });\nanchor.load(\"text\");\nawait context.sync();\n\nif (anchor.items.length > 0) {\n  const insertionPoint = anchor.items[0].getRange(Word.RangeLocation.end);\n  insertionPoint.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Fix typo in BiDS'17 paper: \"Prolusion\" -> \"Propulsion\" in the first-page\n# affiliation footnote (\"* NASA Jet Prolusion Laboratory, California\n# Institute of Technology\"). Word also re-anchors its internal \"last edit\"\n# bookmark (_GoBack) to the location of the correction, so we mirror that\n# too: remove it from its old position and drop it right after the fix.\n\n$d = $word.ActiveDocument\n\n# 1) Locate and fix the misspelled word. It appears exactly once in the\n#    document (a correctly-spelled \"Jet Propulsion Laboratory (JPL)\" exists\n#    elsewhere and must be left untouched), so an exact, whole-word search\n#    is safe and unambiguous.\n$find = $d.Content.Find\n$find.Execute(\"Prolusion\", $true, $true, $false, $false, $false, $true, 1, $false, \"Propulsion\", 2)\n\n# 2) Move the \"_GoBack\" bookmark (Word drops this at the point of the most\n#    recent edit) from wherever it currently sits to right after the word\n#    we just corrected.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n  $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$anchor = $d.Content\n$anchorFind = $anchor.Find\n$anchorFind.Execute(\"NASA Jet Propul\")\nif ($anchorFind.Found) {\n  $anchor.Collapse(0)  # wdCollapseEnd\n  $d.Bookmarks.Add(\"_GoBack\", $anchor)\n}\n"}
